$p = $ppt.ActivePresentation

# --- Slide 2 ("Słabe strony"): merge "Brak " + "logowania" runs into a single run "Brak logowania" ---
$s2 = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(2)
$tr2 = $sh2.TextFrame.TextRange
$para1 = $tr2.Paragraphs(1,1)
# Changing both ends of the paragraph text forces the two runs to collapse into one run;
# then re-assert the exact final text (now a single run, so this just edits its content).
$para1.Text = "XBrak logowaniaY"
$para1.Text = "Brak logowania"

# --- Slide 4 ("Potencjał projektu"): split " Panel" run into " " + "Panel", add "Aspose" paragraph ---
$s4 = $p.Slides.Item(4)
$sh4 = $s4.Shapes.Item(2)
$tr4 = $sh4.TextFrame.TextRange
$para4 = $tr4.Paragraphs(4,1)
# "Admin Panel": "Admin" = chars 1-5, " Panel" = chars 6-11.
# Re-assigning the space character's own text forces a run split at that boundary
# without disturbing "Admin"'s (or "Panel"'s) existing run formatting.
$space = $para4.Characters(6,1)
$space.Text = " "

# Append a brand-new paragraph "Aspose" after the current last paragraph; InsertAfter on the
# whole text range creates a genuine new <a:p> rather than rewriting existing runs.
# (Result assigned to $null so it isn't echoed to the output stream.)
$null = $tr4.InsertAfter("`rAspose")
